$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.675.84'
$ws.Cells.Item(2, 5).Value = '  -4.82%  '
$ws.Cells.Item(3, 4).Value = '3.165.03'
$ws.Cells.Item(3, 5).Value = '  -5.43%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '528.20'
$ws.Cells.Item(5, 5).Value = '  -6.84%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '135.09'
$ws.Cells.Item(6, 5).Value = '  -8.03%  '
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$ws.Cells.Item(8, 4).Value = '3.163.91'
$ws.Cells.Item(8, 5).Value = '  -5.46%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.452'
$ws.Cells.Item(9, 5).Value = '  -6.71%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.20'
$ws.Cells.Item(10, 5).Value = '  -9.23%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.112'
$ws.Cells.Item(11, 5).Value = '  -8.39%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.394'
$ws.Cells.Item(12, 5).Value = '  -5.04%  '
$ws.Cells.Item(13, 4).Value = '3.703.89'
$ws.Cells.Item(13, 5).Value = '  -5.40%  '
$ws.Cells.Item(14, 5).Value = '  -1.57%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '25.86'
$ws.Cells.Item(15, 5).Value = '  -6.74%  '
$ws.Cells.Item(16, 4).Value = '3.158.00'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.0000154'
$ws.Cells.Item(17, 5).Value = '  -8.94%  '
$ws.Cells.Item(18, 4).Value = '57.613.75'
$ws.Cells.Item(18, 5).Value = '  -4.93%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '5.84'
$ws.Cells.Item(19, 5).Value = '  -6.82%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '13.15'
$ws.Cells.Item(20, 5).Value = '  -9.64%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '8.07'
$ws.Cells.Item(21, 5).Value = '  -9.36%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '349.67'
$ws.Cells.Item(22, 5).Value = '  -7.07%  '
$ws.Cells.Item(23, 5).Value = '  +0.06%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '69.63'
$ws.Cells.Item(24, 5).Value = '  -6.91%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.513'
$ws.Cells.Item(25, 5).Value = '  -8.36%  '
$ws.Cells.Item(26, 4).Value = '3.298.19'
$ws.Cells.Item(26, 5).Value = '  -5.50%  '
$ws.Cells.Item(27, 4).Value = '0.0₃0971'
$ws.Cells.Item(27, 5).Value = '  -10.21%  '
$ws.Cells.Item(28, 5).Value = '  -4.35%  '
$ws.Cells.Item(29, 5).Value = '  -0.35%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '6.91'
$ws.Cells.Item(30, 5).Value = '  -5.96%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.13%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.89'
$ws.Cells.Item(32, 5).Value = '  -9.48%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '6.97'
$ws.Cells.Item(33, 5).Value = '  -9.61%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '21.70'
$ws.Cells.Item(34, 5).Value = '  -5.06%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.22'
$ws.Cells.Item(35, 5).Value = '  -6.58%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '4.95'
$ws.Cells.Item(36, 5).Value = '  -6.82%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '158.50'
$ws.Cells.Item(37, 5).Value = '  -5.30%  '
$ws.Cells.Item(38, 5).Value = '  -8.60%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '6.26'
$ws.Cells.Item(39, 5).Value = '  -8.20%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '26.00'
$ws.Cells.Item(40, 5).Value = '  -7.02%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0702'
$ws.Cells.Item(41, 5).Value = '  -6.29%  '
$ws.Cells.Item(42, 4).Value = '3.185.07'
$ws.Cells.Item(42, 5).Value = '  -5.75%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '40.29'
$ws.Cells.Item(43, 5).Value = '  -4.46%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.10'
$ws.Cells.Item(44, 5).Value = '  -3.84%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.698'
$ws.Cells.Item(45, 5).Value = '  -7.87%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.98'
$ws.Cells.Item(46, 5).Value = '  -7.31%  '
$ws.Cells.Item(47, 5).Value = '  -0.09%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.47'
$ws.Cells.Item(48, 5).Value = '  -8.72%  '
$ws.Cells.Item(49, 4).Value = '2.269.09'
$ws.Cells.Item(49, 5).Value = '  -7.57%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '6.23'
$ws.Cells.Item(50, 5).Value = '  -6.65%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '20.73'
$ws.Cells.Item(51, 5).Value = '  -7.29%  '
